$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections to "Nb nouveaux cas positifs" (column C) on a handful
# of existing rows. Column B (cumulative positive cases) is a shared
# formula (IF(TODAY()>A,prev+C,"")) and recalculates automatically. ---
$ws.Range("C241").Value = 756
$ws.Range("C269").Value = 151
$ws.Range("C361").Value = 50
$ws.Range("C363").Value = 89
$ws.Range("C364").Value = 65

# --- Row 365 (2021-02-24) was previously a future/empty placeholder row
# (formulas evaluated to "" since TODAY() <= A364). Fill in the day's
# actual figures. ---
$ws.Range("C365").Value = 10
$ws.Range("E365").Value = 9
$ws.Range("F365").Value = 6
$ws.Range("G365").Value = 26

# L365/M365 are formatted as Text (@). Assigning .Value directly on a
# Text-formatted cell stores the number as a text string (matches real
# Excel COM behaviour). The source rows store these as genuine numbers, so
# temporarily flip the format to General, write the number, then restore
# the original Text format (reusing the existing style, no new styles are
# introduced since both the General and Text variants of this bordered
# style already exist in the workbook).
$ws.Range("L365").NumberFormat = "general"
$ws.Range("M365").NumberFormat = "general"
$ws.Range("L365").Value = 0
$ws.Range("M365").Value = 0
$ws.Range("L365").NumberFormat = "@"
$ws.Range("M365").NumberFormat = "@"
